## Added two extra needy modules for a centurion manual
# - "German Venting Gas" (a translated module, credited to Malde, Tharagon,
#   same contributor/date as the existing "Venting Gas" entry)
# - "Refill That Beer!" (a new needy module by "scripto")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 113: German Venting Gas
$ws.Cells.Item(113, 1).Value = "German Venting Gas"
$ws.Cells.Item(113, 2).Value = "VentGasTranslated"
$ws.Cells.Item(113, 3).Value = 1
$ws.Cells.Item(113, 4).Value = "modules/German Venting Gas.pdf"
$ws.Cells.Item(113, 5).Value = "Malde, Tharagon"
$ws.Cells.Item(113, 6).Value = "2017-01-25"
$ws.Cells.Item(113, 7).Value = 2

# New row 114: Refill That Beer!
$ws.Cells.Item(114, 1).Value = "Refill That Beer!"
$ws.Cells.Item(114, 2).Value = "NeedyBeer"
$ws.Cells.Item(114, 3).Value = 1
$ws.Cells.Item(114, 4).Value = "modules/Refill That Beer!.pdf"
$ws.Cells.Item(114, 5).Value = "scripto"
$ws.Cells.Item(114, 6).Value = "2017-09-23"
$ws.Cells.Item(114, 7).Value = 2

# Restore the sheet view/scroll position & selection as left by the author
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 88
$win.ScrollColumn = 1
$ws.Range("A116").Select() | Out-Null
